$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 7.723979
$ws.Range("H2").Value = 23.171937
$ws.Range("I2").Value = 0.471042132528101
$ws.Range("J2").Value = 0.471042132528101
$ws.Range("M2").Value = 8.142376
$ws.Range("N2").Value = 24.427128
$ws.Range("O2").Value = 0.1741313933276368
$ws.Range("P2").Value = 0.1741313933276368
$ws.Range("Q2").Value = 62.891541234104
$ws.Range("R2").Value = 566.0238711069361
$ws.Range("S2").Value = 0.08202322285313955
$ws.Range("T2").Value = 0.08202322285313955
$ws.Range("G3").Value = 7.723979
$ws.Range("H3").Value = 23.171937
$ws.Range("I3").Value = 0.471042132528101
$ws.Range("J3").Value = 0.471042132528101
$ws.Range("O3").Value = 0.5205382400466131
$ws.Range("P3").Value = 0.5205382400466131
$ws.Range("Q3").Value = 188.0043084834357
$ws.Range("R3").Value = 1692.038776350921
$ws.Range("S3").Value = 0.2451954426539812
$ws.Range("T3").Value = 0.2451954426539812
$ws.Range("G4").Value = 7.723979
$ws.Range("H4").Value = 23.171937
$ws.Range("I4").Value = 0.471042132528101
$ws.Range("J4").Value = 0.471042132528101
$ws.Range("O4").Value = 0.3053303666257501
$ws.Range("P4").Value = 0.3053303666257501
$ws.Range("Q4").Value = 110.277055594086
$ws.Range("R4").Value = 992.493500346774
$ws.Range("S4").Value = 0.1438234670209803
$ws.Range("T4").Value = 0.1438234670209803
$ws.Range("I5").Value = 0.2460132574367717
$ws.Range("J5").Value = 0.2460132574367717
$ws.Range("M5").Value = 8.142376
$ws.Range("N5").Value = 24.427128
$ws.Range("O5").Value = 0.1741313933276368
$ws.Range("P5").Value = 0.1741313933276368
$ws.Range("Q5").Value = 32.84664333778667
$ws.Range("R5").Value = 295.61979004008
$ws.Range("S5").Value = 0.04283863129453565
$ws.Range("T5").Value = 0.04283863129453565
$ws.Range("I6").Value = 0.2460132574367717
$ws.Range("J6").Value = 0.2460132574367717
$ws.Range("O6").Value = 0.5205382400466131
$ws.Range("P6").Value = 0.5205382400466131
$ws.Range("S6").Value = 0.1280593080542715
$ws.Range("T6").Value = 0.1280593080542715
$ws.Range("I7").Value = 0.2460132574367717
$ws.Range("J7").Value = 0.2460132574367717
$ws.Range("O7").Value = 0.3053303666257501
$ws.Range("P7").Value = 0.3053303666257501
$ws.Range("S7").Value = 0.07511531808796455
$ws.Range("T7").Value = 0.07511531808796455
$ws.Range("I8").Value = 0.2829446100351274
$ws.Range("J8").Value = 0.2829446100351274
$ws.Range("M8").Value = 8.142376
$ws.Range("N8").Value = 24.427128
$ws.Range("O8").Value = 0.1741313933276368
$ws.Range("P8").Value = 0.1741313933276368
$ws.Range("Q8").Value = 37.77756039249866
$ws.Range("R8").Value = 339.9980435324881
$ws.Range("S8").Value = 0.04926953917996158
$ws.Range("T8").Value = 0.04926953917996158
$ws.Range("I9").Value = 0.2829446100351274
$ws.Range("J9").Value = 0.2829446100351274
$ws.Range("O9").Value = 0.5205382400466131
$ws.Range("P9").Value = 0.5205382400466131
$ws.Range("S9").Value = 0.1472834893383605
$ws.Range("T9").Value = 0.1472834893383605
$ws.Range("I10").Value = 0.2829446100351274
$ws.Range("J10").Value = 0.2829446100351274
$ws.Range("O10").Value = 0.3053303666257501
$ws.Range("P10").Value = 0.3053303666257501
$ws.Range("S10").Value = 0.08639158151680536
$ws.Range("T10").Value = 0.08639158151680536
